$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format from the existing last header cell (G1) onto the
# new "Save" header cell (H1), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new "Save" column values: 0 for rows 2-13, 1 for row 14
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
$ws.Cells.Item(14, 8).Value = 1
